# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.563.38'
$ws.Range('E2').Value = '  +1.63%  '
$ws.Range('D3').Value = '3.088.75'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '522.68'
$ws.Range('E5').Value = '  +1.28%  '
$ws.Range('D6').Value = '143.23'
$ws.Range('E6').Value = '  +0.72%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +1.25%  '
$ws.Range('D9').Value = '7.34'
$ws.Range('E9').Value = '  +0.47%  '
$ws.Range('E10').Value = '  +0.12%  '
$ws.Range('E11').Value = '  +3.05%  '
$ws.Range('D12').Value = '3.616.01'
$ws.Range('E12').Value = '  -0.20%  '
$ws.Range('E13').Value = '  +0.96%  '
$ws.Range('D14').Value = '26.76'
$ws.Range('E14').Value = '  +4.05%  '
$ws.Range('E15').Value = '  +1.20%  '
$ws.Range('D16').Value = '58.592.05'
$ws.Range('E16').Value = '  +1.52%  '
$ws.Range('D17').Value = '3.086.71'
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('D18').Value = '6.16'
$ws.Range('E18').Value = '  +0.52%  '
$ws.Range('D19').Value = '12.91'
$ws.Range('E19').Value = '  -1.69%  '
$ws.Range('D20').Value = '8.12'
$ws.Range('E20').Value = '  -0.68%  '
$ws.Range('D21').Value = '342.21'
$ws.Range('E21').Value = '  +1.56%  '
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('D24').Value = '65.78'
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('E25').Value = '  +0.22%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').Value = '0.0₃0915'
$ws.Range('E27').Value = '  -1.57%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.60'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.15%  '
$ws.Range('D29').Value = '7.22'
$ws.Range('E29').Value = '  +1.68%  '
$ws.Range('E30').Value = '  +1.64%  '
$ws.Range('D31').Value = '20.99'
$ws.Range('E31').Value = '  +0.12%  '
$ws.Range('E32').Value = '  +2.60%  '
$ws.Range('D33').Value = '154.32'
$ws.Range('E33').Value = '  +0.71%  '
$ws.Range('D34').Value = '4.61'
$ws.Range('E34').Value = '  +1.84%  '
$ws.Range('D35').Value = '6.06'
$ws.Range('E35').Value = '  +2.50%  '
$ws.Range('D36').Value = '26.96'
$ws.Range('E36').Value = '  -1.18%  '
$ws.Range('E37').Value = '  +5.81%  '
$ws.Range('E38').Value = '  -0.91%  '
$ws.Range('D39').Value = '3.129.48'
$ws.Range('E39').Value = '  -0.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.90'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.14%  '
$ws.Range('D41').Value = '36.78'
$ws.Range('E41').Value = '  -0.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.50'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +7.55%  '
$ws.Range('E43').Value = '  -0.71%  '
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('D45').Value = '2.271.78'
$ws.Range('E45').Value = '  -0.90%  '
$ws.Range('E46').Value = '  +1.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '20.80'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.37%  '
$ws.Range('D48').Value = '0.961'
$ws.Range('E48').Value = '  +1.54%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.00'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.70%  '
$ws.Range('B50').Value = 'Bittensor'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D50').Value = '266.31'
$ws.Range('E50').Value = '  +8.40%  '
$ws.Range('B51').Value = 'SuiNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D51').Value = '0.743'
$ws.Range('E51').Value = '  +7.94%  '
